# Fruta / hortaliza, semanal
# Insert one new daily-report row into the "Pepino ensalada" consolidated
# sheet. The new row is inserted at row 787 (pushing the previous rows
# 787-826 down to 788-827) and populated with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 787:826 down to 788:827, creating a blank row 787.
$ws.Rows.Item(787).Insert()

# Populate the newly inserted row with the new record.
$newRow = New-Object 'object[,]' 1,18
$newRow[0,0]  = 5
$newRow[0,1]  = "Macroferia Regional de Talca"
$newRow[0,2]  = "Maule"
$newRow[0,3]  = 45267
$newRow[0,4]  = 7
$newRow[0,5]  = 100112043
$newRow[0,6]  = "Pepino ensalada"
$newRow[0,7]  = "Sin especificar"
$newRow[0,8]  = "Primera"
$newRow[0,9]  = 400
$newRow[0,10] = 16000
$newRow[0,11] = 16000
$newRow[0,12] = 16000
$newRow[0,13] = "`$/caja 80 unidades"
$newRow[0,14] = "Región del Maule"
$newRow[0,15] = 200
$newRow[0,16] = 80
$newRow[0,17] = "Hortaliza"

$ws.Range("A787:R787").Value = $newRow
